# cater no such keyword
$wb = $excel.ActiveWorkbook

$wsCases = $wb.Worksheets.Item("TestCases")
$wsSteps = $wb.Worksheets.Item("TestSteps")

# TestSteps row 2: TestData -> "Browser", ActionKeyword -> "eat" (no such keyword), Results -> FAIL
# (order matters for shared-string allocation: Browser, then eat, then FAIL)
$wsSteps.Range("F2").Value = "Browser"
$wsSteps.Range("E2").Value = "eat"
$wsSteps.Range("G2").Value = "FAIL"

# TestCases!D2 (Results) : PASS -> FAIL
$wsCases.Range("D2").Value = "FAIL"

# Remove the stale Results (PASS) values from rows 3-8, column G
$wsSteps.Range("G3").ClearContents()
$wsSteps.Range("G4").ClearContents()
$wsSteps.Range("G5").ClearContents()
$wsSteps.Range("G6").ClearContents()
$wsSteps.Range("G7").ClearContents()
$wsSteps.Range("G8").ClearContents()

# Selection / active-cell bookkeeping
$null = $wsCases.Range("B13").Select()
$null = $wsSteps.Range("D20").Select()

# TestSteps becomes the active sheet/tab
$null = $wsSteps.Activate()
